$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Experimental (B7): "" -> "true" (literal text, not boolean) ----------
# A literal "true"/"false" typed into a cell is auto-converted by Excel to a
# Boolean. To store it as text instead we enter it with a leading quote
# (forces text) in a scratch cell that is already inside the sheet's used
# range (B13, currently blank) so the sheet <dimension> does not grow, copy
# just the resulting value into the target cell (leaving the target's own
# s="2" style untouched), then clean the scratch cell back to blank with its
# original formatting restored from its s="2" neighbor B14.
$ws.Range("B13").Value = "'true"
$ws.Range("B13").Copy()
$ws.Range("B7").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("B13").ClearContents()
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)  # xlPasteFormats (restore B13's style)

# --- Date (B8): updated timestamp -----------------------------------------
$ws.Range("B8").Value = "2025-07-21T12:46:15+00:00"

# --- Compositional (B18): "" -> "false" (literal text, not boolean) -------
$ws.Range("B13").Value = "'false"
$ws.Range("B13").Copy()
$ws.Range("B18").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("B13").ClearContents()
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)  # xlPasteFormats (restore B13's style)

$excel.CutCopyMode = $false
